$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4500.3335
$ws.Range("I18").Value = 4500.5
$ws.Range("J18").Value = 4500
$ws.Range("K18").Value = 4500.5
$ws.Range("L18").Value = 4500
$ws.Range("M18").Value = -4216.5
$ws.Range("N18").Value = -5068

$ws.Range("H33").Value = 132.5
$ws.Range("I33").Value = 102.77778
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 102.77778
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = 126.22222
$ws.Range("N33").Value = -858

$ws.Range("H53").Value = 520.2
$ws.Range("I53").Value = 520.2
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 520.2
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 116.8

$ws.Range("H62").Value = 2516.875
$ws.Range("I62").Value = 2396
$ws.Range("J62").Value = 2610.889
$ws.Range("K62").Value = 2396
$ws.Range("L62").Value = 2610.889
$ws.Range("M62").Value = -1772
$ws.Range("N62").Value = -3858.889

$ws.Range("H65").Value = 2516.875
$ws.Range("I65").Value = 2396
$ws.Range("J65").Value = 2610.889
$ws.Range("K65").Value = 11980
$ws.Range("L65").Value = 13054.445
$ws.Range("M65").Value = -8860
$ws.Range("N65").Value = -19294.445

$ws.Range("H80").Value = 640
$ws.Range("I80").Value = 757.1429000000001
$ws.Range("J80").Value = 366.66666
$ws.Range("K80").Value = 2271.4287
$ws.Range("L80").Value = 1099.99998
$ws.Range("M80").Value = -1273.4287
$ws.Range("N80").Value = -3095.99998

$ws.Range("H83").Value = 640
$ws.Range("I83").Value = 757.1429000000001
$ws.Range("J83").Value = 366.66666
$ws.Range("K83").Value = 6814.2861
$ws.Range("L83").Value = 3299.99994
$ws.Range("M83").Value = -1822.2861
$ws.Range("N83").Value = -13283.99994

$ws.Range("H103").Value = 752
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 752
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 2256
$ws.Range("N103").Value = -3428
$ws.Range("M103").ClearContents()

$ws.Range("H116").Value = 5750
$ws.Range("I116").Value = 5750
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5750
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -2308

$ws.Range("H132").Value = 2782
$ws.Range("I132").Value = 2914.6667
$ws.Range("J132").Value = 1190
$ws.Range("K132").Value = 8744.000100000001
$ws.Range("L132").Value = 3570
$ws.Range("M132").Value = -6214.000100000001

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()

$ws.Range("H137").Value = 1407.1666
$ws.Range("I137").Value = 1137.7222
$ws.Range("J137").Value = 2215.5
$ws.Range("K137").Value = 3413.1666
$ws.Range("L137").Value = 6646.5
$ws.Range("M137").Value = -863.1665999999996
$ws.Range("N137").Value = -11746.5

$ws.Range("H141").Value = 2598.7693
$ws.Range("I141").Value = 2055.5217
$ws.Range("J141").Value = 6763.6665
$ws.Range("K141").Value = 6166.5651
$ws.Range("L141").Value = 20290.9995
$ws.Range("M141").Value = -986.5650999999998
$ws.Range("N141").Value = -30650.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5226.9688
$ws.Range("I32").Value = 5221.654
$ws.Range("J32").Value = 5250
$ws.Range("K32").Value = 5221.654
$ws.Range("L32").Value = 5250
$ws.Range("M32").Value = -4934.654
$ws.Range("N32").Value = -5824

$ws.Range("H45").Value = 2124.75
$ws.Range("I45").Value = 1999.6666
$ws.Range("J45").Value = 2500
$ws.Range("K45").Value = 1999.6666
$ws.Range("L45").Value = 2500
$ws.Range("M45").Value = -1622.6666
$ws.Range("N45").Value = -3254

$ws.Range("H132").Value = 1749.7142
$ws.Range("I132").Value = 1546.8
$ws.Range("J132").Value = 2257
$ws.Range("K132").Value = 4640.4
$ws.Range("L132").Value = 6771
$ws.Range("M132").Value = -2110.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 995
$ws.Range("I20").Value = 995
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 995
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -748

$ws.Range("H22").Value = 407.69232
$ws.Range("I22").Value = 400
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 400
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -227
$ws.Range("N22").Value = -846

$ws.Range("H105").Value = 1500
$ws.Range("I105").Value = 1500
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1500
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 247

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 70020.60000000001
$ws.Range("I4").Value = 70020.60000000001
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 70020.60000000001
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -69908.60000000001

$ws.Range("H12").Value = 24399.6
$ws.Range("I12").Value = 1000
$ws.Range("J12").Value = 30249.5
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 30249.5
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -30589.5

$ws.Range("H107").Value = 1914.6666
$ws.Range("I107").Value = 1914.6666
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1914.6666
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 5.333399999999983
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 1557.9333
$ws.Range("I132").Value = 1256.1
$ws.Range("J132").Value = 2161.6
$ws.Range("K132").Value = 3768.3
$ws.Range("L132").Value = 6484.799999999999
$ws.Range("M132").Value = -1238.3
$ws.Range("N132").Value = -11544.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34949.668
$ws.Range("I4").Value = 54267.105
$ws.Range("J4").Value = 1583.1818
$ws.Range("K4").Value = 162801.315
$ws.Range("L4").Value = 4749.5454
$ws.Range("M4").Value = -162689.315

$ws.Range("H113").Value = 2400
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 7200
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1444.4445
$ws.Range("I5").Value = 1444.4445
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1444.4445
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1332.4445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1667.2222
$ws.Range("I2").Value = 1.2727273
$ws.Range("J2").Value = 4285.143
$ws.Range("K2").Value = 1.2727273
$ws.Range("L2").Value = 4285.143
$ws.Range("M2").Value = 110.7272727
$ws.Range("N2").Value = -4509.143

$ws.Range("H10").Value = 4112.5
$ws.Range("I10").Value = 2150
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 2150
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -2010
$ws.Range("N10").Value = -10280

$ws.Range("H61").Value = 3562.25
$ws.Range("I61").Value = 2749.6
$ws.Range("J61").Value = 4916.6665
$ws.Range("K61").Value = 2749.6
$ws.Range("L61").Value = 4916.6665
$ws.Range("M61").Value = -2547.6

$ws.Range("H113").Value = 3562.25
$ws.Range("I113").Value = 2749.6
$ws.Range("J113").Value = 4916.6665
$ws.Range("K113").Value = 2749.6
$ws.Range("L113").Value = 4916.6665
$ws.Range("M113").Value = -579.5999999999999

$ws.Range("H122").Value = 2751
$ws.Range("I122").Value = 2751
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8253
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5803

$ws.Range("H132").Value = 7106.0835
$ws.Range("I132").Value = 4628.8335
$ws.Range("J132").Value = 9583.333000000001
$ws.Range("K132").Value = 13886.5005
$ws.Range("L132").Value = 28749.999
$ws.Range("M132").Value = -11356.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8019.2

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H107").Value = 725
$ws.Range("I107").Value = 725
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2175
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -255

$ws.Range("H136").Value = 1298
$ws.Range("I136").Value = 1122.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3367.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -817.5
